$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header label (pandas unnamed-column artifact -> "total")
$ws.Range("B2").Value = "total"

# Remove the empty "situação do domicílio" / "grandes regiões e unidades da
# federação" section-header rows and the trailing "fonte:" footer row, so the
# real data rows close the resulting gaps (values shift up accordingly).
# Deleting from the bottom up keeps the remaining row numbers stable.
$ws.Rows.Item(41).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
